$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing result/profit for two already-settled bets ---
$ws.Range("G35").Value = "Fallo"
$ws.Range("H35").Value = -1

$ws.Range("G55").Value = "Fallo"
$ws.Range("H55").Value = -1

# --- Append new tracker rows (65-75) ---
# Pre-format column B as Text so date-looking strings like "2025-09-02"
# are stored verbatim (matching the rest of the sheet) instead of being
# auto-converted to Excel date serials.
$ws.Range("B65:B75").NumberFormat = "@"

$newRows = @(
    @(14580797, "2025-09-02", "Marco Cecchinato", "Dmitry Popko", "Gana Dmitry Popko", 2.5),
    @(14592745, "2025-09-02", "Nicolai Budkov Kjaer", "Gianluca Cadenasso", "Gana Gianluca Cadenasso", 4.5),
    @(14579391, "2025-09-02", "Giulio Zeppieri", "Gauthier Onclin", "Gana Gauthier Onclin", 2.38),
    @(14592265, "2025-09-02", "Filip Peliwo", "Luca Van Assche", "Gana Filip Peliwo", 5),
    @(14592262, "2025-09-02", "Kasidit Samrej", "Tung-Lin Wu", "Gana Kasidit Samrej", 2.75),
    @(14592263, "2025-09-02", "Yan Bai", "Tsung-Hao Huang", "Gana Tsung-Hao Huang", 2.2),
    @(14581493, "2025-09-02", "Oleg Prihodko", "Lukas Neumayer", "Gana Oleg Prihodko", 3.4),
    @(14591061, "2025-09-02", "Rudolf Molleker", "David Poljak", "Gana David Poljak", 8),
    @(14591059, "2025-09-02", "Kai Wehnelt", "Neil Oberleitner", "Gana Kai Wehnelt", 2.75),
    @(14591245, "2025-09-02", "Michal Krajci", "Milos Karol", "Gana Michal Krajci", 2.75),
    @(14591246, "2025-09-02", "Peter Makk", "Calvin Hemery", "Gana Peter Makk", 3)
)

$r = 65
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Restore the default "Normal" style on column B so the new cells don't
# carry a leftover explicit number-format style (matches original sheet
# where date cells have no style attribute).
$ws.Range("B65:B75").Style = "Normal"
